$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.089.81"
$ws.Range("E2").Value = "'  +0.07%  "
$ws.Range("D3").Value = "'1.821.70"
$ws.Range("E3").Value = "'  -0.52%  "
$ws.Range("D4").Value = "'0.9986"
$ws.Range("E4").Value = "'  -0.17%  "
$ws.Range("D5").Value = "'241.24"
$ws.Range("E5").Value = "'  -0.85%  "
$ws.Range("D6").Value = "'0.6153"
$ws.Range("E6").Value = "'  -2.06%  "
$ws.Range("D7").Value = "'0.9989"
$ws.Range("E7").Value = "'  -0.17%  "
$ws.Range("D8").Value = "'0.07328"
$ws.Range("E8").Value = "'  -2.14%  "
$ws.Range("D9").Value = "'0.2890"
$ws.Range("E9").Value = "'  -1.05%  "
$ws.Range("D10").Value = "'22.91"
$ws.Range("E10").Value = "'  -1.40%  "
$ws.Range("D11").Value = "'0.07657"
$ws.Range("E11").Value = "'  -0.23%  "
$ws.Range("D12").Value = "'1.820.85"
$ws.Range("E12").Value = "'  -0.91%  "
$ws.Range("D13").Value = "'4.943"
$ws.Range("E13").Value = "'  -1.19%  "
$ws.Range("D14").Value = "'0.6591"
$ws.Range("E14").Value = "'  -1.21%  "
$ws.Range("D15").Value = "'81.69"
$ws.Range("E15").Value = "'  -1.19%  "
$ws.Range("D16").Value = "'0.000008975"
$ws.Range("E16").Value = "'  -3.86%  "
$ws.Range("D17").Value = "'5.824"
$ws.Range("E17").Value = "'  -2.59%  "
$ws.Range("D18").Value = "'29.048.79"
$ws.Range("E18").Value = "'  -0.12%  "
$ws.Range("D19").Value = "'2.061.26"
$ws.Range("E19").Value = "'  -0.85%  "
$ws.Range("D20").Value = "'237.97"
$ws.Range("E20").Value = "'  +6.72%  "
$ws.Range("E21").Value = "'  -1.31%  "
$ws.Range("D22").Value = "'0.9986"
$ws.Range("E22").Value = "'  -0.35%  "
$ws.Range("D23").Value = "'7.104"
$ws.Range("E23").Value = "'  +0.13%  "
$ws.Range("E24").Value = "'  -0.15%  "
$ws.Range("D25").Value = "'157.27"
$ws.Range("E25").Value = "'  -1.65%  "
$ws.Range("D26").Value = "'0.1406"
$ws.Range("E26").Value = "'  +0.94%  "
$ws.Range("D27").Value = "'8.402"
$ws.Range("E27").Value = "'  -0.94%  "
$ws.Range("D28").Value = "'17.59"
$ws.Range("E28").Value = "'  -1.49%  "
$ws.Range("D29").Value = "'1.481"
$ws.Range("E29").Value = "'  -1.13%  "
$ws.Range("D30").Value = "'0.05550"
$ws.Range("E30").Value = "'  -2.22%  "
$ws.Range("E31").Value = "'  +0.42%  "
$ws.Range("D32").Value = "'4.088"
$ws.Range("E32").Value = "'  -1.40%  "
$ws.Range("D33").Value = "'1.206"
$ws.Range("E33").Value = "'  -0.06%  "
$ws.Range("D34").Value = "'0.7336"
$ws.Range("E34").Value = "'  -1.06%  "
$ws.Range("D35").Value = "'1.814"
$ws.Range("E35").Value = "'  -0.92%  "
$ws.Range("D36").Value = "'1.130"
$ws.Range("E36").Value = "'  -0.86%  "
$ws.Range("D37").Value = "'2.604"
$ws.Range("E37").Value = "'  -2.38%  "
$ws.Range("D38").Value = "'2.828"
$ws.Range("E38").Value = "'  +2.53%  "
$ws.Range("D39").Value = "'1.208.28"
$ws.Range("E39").Value = "'  -0.53%  "
$ws.Range("D40").Value = "'0.01754"
$ws.Range("E40").Value = "'  -1.21%  "
$ws.Range("D41").Value = "'6.348"
$ws.Range("E41").Value = "'  -2.39%  "
$ws.Range("D42").Value = "'0.8926"
$ws.Range("E42").Value = "'  +0.42%  "
$ws.Range("D43").Value = "'0.9985"
$ws.Range("E43").Value = "'  -0.19%  "
$ws.Range("D44").Value = "'100.84"
$ws.Range("E44").Value = "'  -1.15%  "
$ws.Range("D45").Value = "'1.970.72"
$ws.Range("E45").Value = "'  -0.47%  "
$ws.Range("D46").Value = "'64.55"
$ws.Range("E46").Value = "'  -1.59%  "
$ws.Range("D47").Value = "'0.5082"
$ws.Range("E47").Value = "'  -0.13%  "
$ws.Range("E48").Value = "'  -5.65%  "
$ws.Range("D49").Value = "'0.3995"
$ws.Range("E49").Value = "'  -1.62%  "
$ws.Range("D50").Value = "'9.052"
$ws.Range("E50").Value = "'  +0.75%  "
$ws.Range("D51").Value = "'0.05750"
$ws.Range("E51").Value = "'  -1.18%  "
